# Update countries & provincias Spain
# Applies the data refresh captured in the target diff:
#  - bumps the "Datos actualizados" timestamp from 07:05 to 08:05
#  - refreshes case counters for several countries (India, Israel,
#    Afganistan, Australia, Uzbekistan, Taiwan)
#  - refreshes counters for El Salvador / Bulgaria / Bosnia y Herzegovina
#    and re-sorts those three rows (El Salvador, Bulgaria, Bosnia y
#    Herzegovina) since El Salvador's totals overtook Bulgaria's
#  - swaps the row order of Belice / Santa Lucia (same totals, re-sorted)
#  - swaps the row order of Islas Virgenes Britanicas / Papua Nueva
#    Guinea (same totals, re-sorted)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp row ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 31 de Mayo de 2020 a las 08:05"

# --- Row 12: India -------------------------------------------------------
$ws.Range("B12").Value = 182490
$ws.Range("C12").Value = 663
$ws.Range("E12").Value = 90320
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 5186

# --- Row 43: Israel --------------------------------------------------------
$ws.Range("B43").Value = 17024
$ws.Range("C43").Value = 12
$ws.Range("D43").Value = 14812
$ws.Range("E43").Value = 1928

# --- Row 48: Afganistan ---------------------------------------------------
$ws.Range("B48").Value = 15205
$ws.Range("C48").Value = 680
$ws.Range("D48").Value = 1328
$ws.Range("E48").Value = 13620
$ws.Range("G48").Value = 8
$ws.Range("H48").Value = 257

# --- Row 66: Australia -----------------------------------------------------
$ws.Range("B66").Value = 7195
$ws.Range("C66").Value = 10
$ws.Range("E66").Value = 478

# --- Row 78: Uzbekistan ------------------------------------------------------
$ws.Range("B78").Value = 3554
$ws.Range("C78").Value = 8
$ws.Range("E78").Value = 757

# --- Rows 86-88: El Salvador / Bulgaria / Bosnia y Herzegovina re-sort ------
# New order: El Salvador (86), Bulgaria (87), Bosnia y Herzegovina (88)
$ws.Range("A86").Value = "El Salvador"
$ws.Range("B86").Value = 2517
$ws.Range("C86").Value = 122
$ws.Range("D86").Value = 1040
$ws.Range("E86").Value = 1431
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 2
$ws.Range("H86").Value = 46

$ws.Range("A87").Value = "Bulgaria"
$ws.Range("B87").Value = 2513
$ws.Range("C87").Value = 14
$ws.Range("D87").Value = 1074
$ws.Range("E87").Value = 1299
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 140

$ws.Range("A88").Value = "Bosnia y Herzegovina"
$ws.Range("B88").Value = 2494
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 1831
$ws.Range("E88").Value = 510
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 153

# --- Row 143: Taiwan ----------------------------------------------------
$ws.Range("D143").Value = 423
$ws.Range("E143").Value = 12

# --- Rows 200-201: Belice / Santa Lucia re-sort -----------------------------
$ws.Range("A200").Value = "Belice"
$ws.Range("B200").Value = 18
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 16
$ws.Range("E200").Value = 0
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 2

$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("B201").Value = 18
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 18
$ws.Range("E201").Value = 0
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 0

# --- Rows 213-214: Islas Virgenes Britanicas / Papua Nueva Guinea re-sort ---
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("B213").Value = 8
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 7
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("B214").Value = 8
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 8
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0
